$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 1517
$ws1.Range("F5").Value = 7760
$ws1.Range("F6").Value = 4868
$ws1.Range("F7").Value = 7157
$ws1.Range("F26").Value = 159
$ws1.Range("F28").Value = 12
$ws1.Range("F32").Value = 2

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F8").Value = 99
$ws3.Range("F9").Value = 1721
$ws3.Range("F10").Value = 2627

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 1517
$ws4.Range("F8").Value = 7760
$ws4.Range("F10").Value = 4868
$ws4.Range("F11").Value = 7157
$ws4.Range("F15").Value = 99
$ws4.Range("F17").Value = 1721
$ws4.Range("F18").Value = 2627
$ws4.Range("F28").Value = 159
